# Applies two changes:
#  1. Swaps the table style on the B1/B2 table (slide 5) from the
#     "Integral"-deck default table style to the new style id.
#  2. Re-colours the deck's theme (used by the slide master, i.e.
#     ppt/theme/theme1.xml) from the "Integral" / "Red Violet" palette
#     over to the stock "Office" palette, matching the colour values
#     that the target theme (formerly the notes-master-only theme)
#     uses.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $s = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $s.Shapes.Count; $shi++) {
        $sh = $s.Shapes.Item($shi)
        if ($sh.HasTable) {
            $sh.Table.ApplyStyle("{D9232DCC-7D64-4743-A861-C6FA58EFBF4E}")
        }
    }
}

# --- 2. Theme colour scheme ------------------------------------------
# Office RGB() encodings (0x00BBGGRR) for the 12 standard theme colours.
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
